$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 264.33334
$ws.Range("I4").Value = 264.33334
$ws.Range("K4").Value = 264.33334
$ws.Range("M4").Value = -150.33334
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H107").Value = 670
$ws.Range("I107").Value = 670
$ws.Range("K107").Value = 670
$ws.Range("M107").Value = 1250
$ws.Range("H135").Value = 1835
$ws.Range("I135").Value = 1253
$ws.Range("J135").Value = 2999
$ws.Range("K135").Value = 11277
$ws.Range("L135").Value = 26991
$ws.Range("M135").Value = -8742
$ws.Range("N135").Value = -32061
$ws.Range("H137").Value = 1806.5
$ws.Range("I137").Value = 1575.7778
$ws.Range("J137").Value = 3883
$ws.Range("K137").Value = 4727.3334
$ws.Range("L137").Value = 11649
$ws.Range("M137").Value = -2177.3334
$ws.Range("N137").Value = -16749
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2098.2727
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2226
$ws.Range("H88").Value = 2372.2104
$ws.Range("J88").Value = 2974.5
$ws.Range("L88").Value = 2974.5
$ws.Range("N88").Value = -3786.5
$ws.Range("H91").Value = 2372.2104
$ws.Range("J91").Value = 2974.5
$ws.Range("L91").Value = 2974.5
$ws.Range("N91").Value = -5782.5
$ws.Range("H116").Value = 2098.2727
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("N116").Value = -6588
$ws.Range("H132").Value = 2484
$ws.Range("I132").Value = 2484
$ws.Range("K132").Value = 7452
$ws.Range("M132").Value = -4922
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2098.2727
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2228
$ws.Range("H105").Value = 1912.3
$ws.Range("I105").Value = 1953.5
$ws.Range("K105").Value = 1953.5
$ws.Range("M105").Value = -206.5
$ws.Range("H134").Value = 9156.267
$ws.Range("I134").Value = 9528.75
$ws.Range("K134").Value = 28586.25
$ws.Range("M134").Value = -26051.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4406.25
$ws.Range("I31").Value = 4305.5
$ws.Range("J31").Value = 4507
$ws.Range("K31").Value = 4305.5
$ws.Range("L31").Value = 4507
$ws.Range("M31").Value = -4010.5
$ws.Range("N31").Value = -5097
$ws.Range("H34").Value = 4406.25
$ws.Range("I34").Value = 4305.5
$ws.Range("J34").Value = 4507
$ws.Range("K34").Value = 4305.5
$ws.Range("L34").Value = 4507
$ws.Range("M34").Value = -4103.5
$ws.Range("N34").Value = -4911
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 22500
$ws.Range("I55").Value = 22500
$ws.Range("K55").Value = 22500
$ws.Range("M55").Value = -22185
$ws.Range("H132").Value = 1240.6666
$ws.Range("I132").Value = 1152.6666
$ws.Range("K132").Value = 3457.9998
$ws.Range("M132").Value = -927.9998000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 9007.25
$ws.Range("I81").Value = 6000
$ws.Range("K81").Value = 18000
$ws.Range("M81").Value = -16877
$ws.Range("H84").Value = 9007.25
$ws.Range("I84").Value = 6000
$ws.Range("K84").Value = 54000
$ws.Range("M84").Value = -48384
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H80").Value = 1983.25
$ws.Range("I80").Value = 1774.2858
$ws.Range("J80").Value = 2041.76
$ws.Range("K80").Value = 1774.2858
$ws.Range("L80").Value = 2041.76
$ws.Range("M80").Value = -776.2858000000001
$ws.Range("N80").Value = -4037.76
$ws.Range("H83").Value = 1983.25
$ws.Range("I83").Value = 1774.2858
$ws.Range("J83").Value = 2041.76
$ws.Range("K83").Value = 8871.429
$ws.Range("L83").Value = 10208.8
$ws.Range("M83").Value = -3879.429
$ws.Range("N83").Value = -20192.8
$ws.Range("H132").Value = 1700.375
$ws.Range("I132").Value = 1700.375
$ws.Range("K132").Value = 5101.125
$ws.Range("M132").Value = -2571.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H40").Value = 7723.625
$ws.Range("I40").Value = 7798.1665
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 7798.1665
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -7662.1665
$ws.Range("N40").Value = -7772
$ws.Range("H43").Value = 10011.2
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10011.2
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10011.2
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10397.2
$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 500
$ws.Range("L68").Value = 500
$ws.Range("M68").Value = 249
$ws.Range("N68").Value = -1998
$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 2500
$ws.Range("L71").Value = 2500
$ws.Range("M71").Value = 1244
$ws.Range("N71").Value = -9988
$ws.Range("H82").Value = 1426.25
$ws.Range("I82").Value = 852.6667
$ws.Range("K82").Value = 852.6667
$ws.Range("M82").Value = -491.6667
$ws.Range("H85").Value = 1426.25
$ws.Range("I85").Value = 852.6667
$ws.Range("K85").Value = 852.6667
$ws.Range("M85").Value = 395.3333
$ws.Range("H100").Value = 1998
$ws.Range("I100").Value = 1998
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1998
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1457
$ws.Range("N100").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2444.389
$ws.Range("J2").Value = 1999
$ws.Range("L2").Value = 1999
$ws.Range("N2").Value = -2223
$ws.Range("H4").Value = 201666.67
$ws.Range("I4").Value = 201666.67
$ws.Range("K4").Value = 201666.67
$ws.Range("M4").Value = -201553.67
$ws.Range("H6").Value = 78
$ws.Range("J6").Value = 78
$ws.Range("L6").Value = 78
$ws.Range("N6").Value = -308
$ws.Range("H68").Value = 24665
$ws.Range("J68").Value = 24665
$ws.Range("L68").Value = 24665
$ws.Range("N68").Value = -26287
$ws.Range("H71").Value = 24665
$ws.Range("J71").Value = 24665
$ws.Range("L71").Value = 73995
$ws.Range("N71").Value = -82107
$ws.Range("H81").Value = 400
$ws.Range("J81").Value = 600
$ws.Range("L81").Value = 1200
$ws.Range("N81").Value = -3322
$ws.Range("H84").Value = 400
$ws.Range("J84").Value = 600
$ws.Range("L84").Value = 6000
$ws.Range("N84").Value = -16608
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 2507.3333
$ws.Range("I122").Value = 2507.3333
$ws.Range("K122").Value = 7521.999899999999
$ws.Range("M122").Value = -5071.999899999999
$ws.Range("H136").Value = 1733
$ws.Range("I136").Value = 1644.65
$ws.Range("K136").Value = 4933.950000000001
$ws.Range("M136").Value = -2383.950000000001
